# Fill in the title page / abstract / table-of-contents placeholder content
# for the EE310 Team 16 report.  The document originally contains a single
# empty paragraph; we keep that paragraph in place (it becomes the blank
# line after the title) and insert the rest of the front-matter paragraphs
# around it using raw WordprocessingML via Range.InsertXML, which lets us
# faithfully reproduce details (like the w:proofErr spell-check wrappers
# and explicit page breaks) that a plain TypeText pass would not add.

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument

# --- Content that goes AFTER the (pre-existing) blank paragraph ---------
$afterXml = @"
<w:p $ns><w:proofErr w:type="spellStart"/><w:r><w:t>Aliaxis</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
<w:p $ns/>
<w:p $ns><w:r><w:t>Authors:</w:t></w:r></w:p>
<w:p $ns><w:proofErr w:type="spellStart"/><w:r><w:t>Yunbo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Hu</w:t></w:r></w:p>
<w:p $ns><w:r><w:t>Damian Lee</w:t></w:r></w:p>
<w:p $ns><w:proofErr w:type="spellStart"/><w:r><w:t>Hawon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ryu</w:t></w:r></w:p>
<w:p $ns><w:r><w:t>Aldonza Watt</w:t></w:r></w:p>
<w:p $ns><w:r><w:t>Team 16</w:t></w:r></w:p>
<w:p $ns><w:r><w:t>Date: 23 April</w:t></w:r></w:p>
<w:p $ns><w:r><w:br w:type="page"/></w:r></w:p>
<w:p $ns><w:r><w:lastRenderedPageBreak/><w:t>Abstract</w:t></w:r></w:p>
<w:p $ns/>
<w:p $ns><w:r><w:br w:type="page"/></w:r></w:p>
<w:p $ns><w:r><w:lastRenderedPageBreak/><w:t>Table of contents</w:t></w:r></w:p>
<w:p $ns/>
<w:p $ns><w:r><w:br w:type="page"/></w:r></w:p>
"@

# --- Content that goes BEFORE the (pre-existing) blank paragraph --------
$beforeXml = "<w:p $ns><w:r><w:t>Water Level Detector PoC</w:t></w:r></w:p>"

# The document starts out with exactly one (empty) paragraph. Grab it
# first, since indices/ranges shift once we start inserting.
$firstPara = $d.Paragraphs.Item(1)

# Create a fresh empty paragraph right after it, then pour the "after"
# WordML into that new paragraph (InsertXML on a collapsed range behaves
# like pasting at the cursor, splitting the target paragraph as needed).
[void]$firstPara.Range.InsertParagraphAfter()
$afterHost = $d.Paragraphs.Item(2)
[void]$afterHost.Range.InsertXML($afterXml)

# Finally, insert the title paragraph at the very start of the document,
# ahead of the original blank paragraph.
$startRange = $d.Range(0, 0)
[void]$startRange.InsertXML($beforeXml)
